$d = $word.ActiveDocument

$d.Content.Find.Execute("return a vector of correlations", $true, $false, $false, $false, $false,
                         $true, 1, $false, "for i", 2)
